$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8594.505859375,1,0,0),
    @(1900.897094726562,2,2.222222222222222,2.222222222222222),
    @(1653.625610351562,2,2.857142857142857,5.079365079365079),
    @(1653.625610351562,2,4.444444444444445,9.523809523809524),
    @(4524.21337890625,4,3.333333333333333,12.85714285714286),
    @(8478.298828125,6,2,14.85714285714286),
    @(18554.083984375,10,3.75,18.60714285714286),
    @(41376.64453125,17,5.185185185185185,23.79232804232804),
    @(50227.77734375,20,2.162162162162162,25.95449020449021),
    @(50799.57421875,21,0.975609756097561,26.93009996058777),
    @(24190.767578125,18,2.051282051282051,28.98138201186982),
    @(4393.96337890625,9,6.108304630544683,35.0896866424145),
    @(14.30727481842041,3,2.357022603955159,37.44670924636966),
    @(16694.365234375,9,6.009252125773315,43.45596137214297),
    @(14.30727481842041,2,1.818181818181818,45.27414319032479),
    @(2511.444091796875,3,8.944271909999159,54.21841510032395),
    @(22004.236328125,11,8.081220356417687,62.29963545674164),
    @(8865.7099609375,10.63749142467192,9.103511685801108,71.40314714254275),
    @(5179.53369140625,9,2.880639959127004,74.28378710166976),
    @(26740.583984375,12,2.1295885499998,76.41337565166955),
    @(14.30727481842041,6,1.111111111111111,77.52448676278067),
    @(15735.3154296875,11.01196482921286,10.04469951718429,87.56918627996497),
    @(34399.63671875,14,1.787998658056773,89.35718493802175),
    @(37390.48046875,16,1.49071198499986,90.84789692302161),
    @(16321.2607421875,14,3.590109871423003,94.4380067944446),
    @(14.30727481842041,6,1,95.4380067944446),
    @(14554.37890625,10,6.73145600891813,102.1694628033627),
    @(5053.35400390625,9,1.488645855129574,103.6581086584923),
    @(2468.705810546875,5,5.890150893739516,109.5482595522318),
    @(12573.4228515625,7,1.666666666666667,111.2149262188985),
    @(8779.3212890625,8,1.333333333333333,112.5482595522318),
    @(35344.87890625,16,8.158431221748456,120.7066907739803),
    @(6705.90087890625,11.66170647106067,5.829183211623753,126.535873985604),
    @(51982.15234375,20,5.649898825367719,132.1857728109718),
    @(31559.291015625,18,5.668594533825794,137.8543673447976),
    @(14.30727481842041,2,4.123105625617661,141.9774729704152),
    @(1644.467407226562,2,7.071067811865476,149.0485407822807),
    @(5199.88525390625,6,7.905694150420948,156.9542349327016),
    @(11461.6923828125,9,4,160.9542349327016),
    @(36721.51171875,16,7.244308110509934,168.1985430432116),
    @(63091.8984375,21,2.702702702702703,170.9012457459143),
    @(14.30727481842041,5,1.720052290384454,172.6212980362987),
    @(8538.3310546875,7,2.357022603955159,174.9783206402539),
    @(437.1109924316406,1.917078252831711,9.247660519987756,184.2259811602416),
    @(4695.16650390625,4,3.380046561058861,187.6060277213005),
    @(4247.43212890625,5.338293528939332,4.789029110233523,192.395056831534),
    @(16421.029296875,11,10.53026162347456,202.9253184550086),
    @(21095.669921875,12,1.229750923802691,204.1550693788113),
    @(6351.78369140625,10.8799729114185,1.954607189577446,206.1096765683887),
    @(18366.142578125,12.78491521734411,8.823457311189456,214.9331338795782),
    @(36686.33984375,17,5.413651635042291,220.3467855146205),
    @(50221.10546875,19,1.111111111111111,221.4578966257316),
    @(54190.05859375,20,0.5128205128205128,221.9707171385521),
    @(61460.26171875,22,0.9523809523809523,222.9230980909331),
    @(37622.54296875,20,4.285714285714286,227.2088123766474),
    @(5777.86962890625,11,6.735681618651967,233.9444939952993),
    @(5692.39306640625,9.541242364689118,3.51054839960607,237.4550423949054),
    @(5737.16650390625,9,5.393381847510125,242.8484242424155),
    @(30094.345703125,14,4.682752006203916,247.5311762486194),
    @(25535.259765625,15,4.876598490941707,252.4077747395611),
    @(39161.56640625,16,0.6451612903225806,253.0529360298837),
    @(46474.91015625,18,1.176470588235294,254.229406618119),
    @(53975.96484375,20,1.052631578947368,255.2820381970664),
    @(65647.4609375,25,4,259.2820381970664),
    @(12452.5263671875,13,11.57894736842105,270.8609855654875),
    @(12876.4794921875,12.43989448617086,2.4860776540464,273.3470632195339),
    @(27978.830078125,15,5.153859330930844,278.5009225504647),
    @(43222.87109375,18,2.424242424242424,280.9251649747072),
    @(6596.00244140625,12,5.497474167490214,286.4226391421974)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
